$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 66 ("Femacal de La Calera" /
# Papaya weekly block). This shifts the existing rows 66-76 down to 68-78,
# carrying their values/styles with them (matches the diff, which is a
# down-shift of the whole historical block to make room for a new week).
$ws.Rows("66:67").Insert()

# Row 66: new "Primera" quality record for the week of 2023-07-24 (serial 45131)
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(66, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = 45131
$ws.Cells.Item(66, 5).Value = 5
$ws.Cells.Item(66, 6).Value = "Fruta"
$ws.Cells.Item(66, 7).Value = 100108
$ws.Cells.Item(66, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(66, 9).Value = 100108004
$ws.Cells.Item(66, 10).Value = "Papaya"
$ws.Cells.Item(66, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(66, 12).Value = "Primera"
$ws.Cells.Item(66, 13).Value = 56
$ws.Cells.Item(66, 14).Value = 20000
$ws.Cells.Item(66, 15).Value = 20000
$ws.Cells.Item(66, 16).Value = 20000
$ws.Cells.Item(66, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(66, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(66, 19).Value = 2000
$ws.Cells.Item(66, 20).Value = 10

# Row 67: new "Segunda" quality record for the same week
$ws.Cells.Item(67, 1).Value = 3
$ws.Cells.Item(67, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(67, 3).Value = "Coquimbo"
$ws.Cells.Item(67, 4).Value = 45131
$ws.Cells.Item(67, 5).Value = 5
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100108
$ws.Cells.Item(67, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(67, 9).Value = 100108004
$ws.Cells.Item(67, 10).Value = "Papaya"
$ws.Cells.Item(67, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(67, 12).Value = "Segunda"
$ws.Cells.Item(67, 13).Value = 50
$ws.Cells.Item(67, 14).Value = 17000
$ws.Cells.Item(67, 15).Value = 17000
$ws.Cells.Item(67, 16).Value = 17000
$ws.Cells.Item(67, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(67, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(67, 19).Value = 1700
$ws.Cells.Item(67, 20).Value = 10
